$wb = $excel.ActiveWorkbook

# Sheet "OFF" - update row 2 (H) values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 528
$wsOff.Range("C2").Value = 362
$wsOff.Range("D2").Value = 122
$wsOff.Range("E2").Value = 56

# Sheet "DEF" - update row 2 (H) values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 485
$wsDef.Range("C2").Value = 335
$wsDef.Range("D2").Value = 111
$wsDef.Range("F2").Value = 12
